$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3716865831074188
$ws.Range("C2").Value = 0.05611706257801075
$ws.Range("D2").Value = 0.3069392310650585
$ws.Range("F2").Value = 1.382777874846866
$ws.Range("G2").Value = 0.002452010522008441
$ws.Range("I2").Value = 0.6277365463388946
$ws.Range("J2").Value = 0.349601523791975
$ws.Range("K2").Value = 0.3997780441873999
$ws.Range("N2").Value = 1.488755955518869
$ws.Range("O2").Value = 3.072261719399762
$ws.Range("B3").Value = 0.3350010994907393
$ws.Range("C3").Value = 0.04955965120919359
$ws.Range("D3").Value = 0.2961529380605583
$ws.Range("F3").Value = 1.382521390128545
$ws.Range("G3").Value = 0.002454346493860032
$ws.Range("I3").Value = 0.6330069154511762
$ws.Range("J3").Value = 0.3382217178279774
$ws.Range("K3").Value = 0.3585184591513553
$ws.Range("N3").Value = 1.504171835118195
$ws.Range("O3").Value = 3.086941397435169
$ws.Range("B4").Value = 0.3125228076359008
$ws.Range("C4").Value = 0.04553123094453326
$ws.Range("D4").Value = 0.2896702114335739
$ws.Range("F4").Value = 1.383101361320662
$ws.Range("G4").Value = 0.002455857848092324
$ws.Range("I4").Value = 0.6365569477843316
$ws.Range("J4").Value = 0.3314377112814242
$ws.Range("K4").Value = 0.33322143202912
$ws.Range("N4").Value = 1.514132284327914
$ws.Range("O4").Value = 3.097574637028231
$ws.Range("B5").Value = 0.3033750314684482
$ws.Range("C5").Value = 0.04388914327655868
$ws.Range("D5").Value = 0.2870638096966616
$ws.Range("F5").Value = 1.383523216359123
$ws.Range("G5").Value = 0.002456493171498911
$ws.Range("I5").Value = 0.6380825670879098
$ws.Range("J5").Value = 0.3287242745443706
$ws.Range("K5").Value = 0.3229224487375859
$ws.Range("N5").Value = 1.518315740423141
$ws.Range("O5").Value = 3.102315106829764
$ws.Range("B6").Value = 0.3018568100973482
$ws.Range("C6").Value = 0.04361644918343188
$ws.Range("D6").Value = 0.2866331583527142
$ws.Range("F6").Value = 1.383604471717767
$ws.Range("G6").Value = 0.002456599841832462
$ws.Range("I6").Value = 0.6383406633488526
$ws.Range("J6").Value = 0.3282767966717302
$ws.Range("K6").Value = 0.3212129168984745
$ws.Range("N6").Value = 1.519017919514104
$ws.Range("O6").Value = 3.103126860863668
$ws.Range("B7").Value = 0.3123993869457991
$ws.Range("C7").Value = 0.04550908700790046
$ws.Range("D7").Value = 0.2896349172330162
$ws.Range("F7").Value = 1.383106299377921
$ws.Range("G7").Value = 0.002455866337605897
$ws.Range("I7").Value = 0.6365772031541788
$ws.Range("J7").Value = 0.331400910024982
$ws.Range("K7").Value = 0.333082496054999
$ws.Range("N7").Value = 1.514188199800891
$ws.Range("O7").Value = 3.097636919351785
$ws.Range("B8").Value = 0.359028048957498
$ws.Range("C8").Value = 0.05385654908648974
$ws.Range("D8").Value = 0.3031911149150375
$ws.Range("F8").Value = 1.382536420414247
$ws.Range("G8").Value = 0.002452800006159315
$ws.Range("I8").Value = 0.6294885938842611
$ws.Range("J8").Value = 0.3456355766519579
$ws.Range("K8").Value = 0.3855444913570807
$ws.Range("N8").Value = 1.493968613805382
$ws.Range("O8").Value = 3.076987073727523
$ws.Range("B9").Value = 0.4508176999664215
$ws.Range("C9").Value = 0.07020682182916005
$ws.Range("D9").Value = 0.330882221567748
$ws.Range("F9").Value = 1.387268866178218
$ws.Range("G9").Value = 0.002447395750386505
$ws.Range("I9").Value = 0.6180797620536467
$ws.Range("J9").Value = 0.3751644619220258
$ws.Range("K9").Value = 0.4886918730699392
$ws.Range("N9").Value = 1.458244915103419
$ws.Range("O9").Value = 3.04934732537501
$ws.Range("B10").Value = 0.5184494680599414
$ws.Range("C10").Value = 0.08220613077361349
$ws.Range("D10").Value = 0.3518987739395811
$ws.Range("F10").Value = 1.39431294574328
$ws.Range("G10").Value = 0.002443792730440398
$ws.Range("I10").Value = 0.6112176000598595
$ws.Range("J10").Value = 0.3978493844656725
$ws.Range("K10").Value = 0.5646188681002968
$ws.Range("N10").Value = 1.434390394496843
$ws.Range("O10").Value = 3.036882433331755
$ws.Range("B11").Value = 0.5492551599248543
$ws.Range("C11").Value = 0.0876617824243624
$ws.Range("D11").Value = 0.3616050986888695
$ws.Range("F11").Value = 1.398292423742674
$ws.Range("G11").Value = 0.00244223264122174
$ws.Range("I11").Value = 0.6084260499317971
$ws.Range("J11").Value = 0.4083856902826142
$ws.Range("K11").Value = 0.5991876867117014
$ws.Range("N11").Value = 1.424057227670025
$ws.Range("O11").Value = 3.03291610736224
$ws.Range("B12").Value = 0.5609256989036737
$ws.Range("C12").Value = 0.08972722928859866
$ws.Range("D12").Value = 0.3653014943049584
$ws.Range("F12").Value = 1.399910785216704
$ws.Range("G12").Value = 0.002441653168974131
$ws.Range("I12").Value = 0.6074164454708253
$ws.Range("J12").Value = 0.4124067426155733
$ws.Range("K12").Value = 0.612281717026093
$ws.Range("N12").Value = 1.420218866854992
$ws.Range("O12").Value = 3.031659272876254
$ws.Range("B13").Value = 0.5584120227841254
$ws.Range("C13").Value = 0.0892824213201493
$ws.Range("D13").Value = 0.3645044856365587
$ws.Range("F13").Value = 1.399557287307147
$ws.Range("G13").Value = 0.002441777467080484
$ws.Range("I13").Value = 0.6076317690301281
$ws.Range("J13").Value = 0.411539350179126
$ws.Range("K13").Value = 0.6094615347269041
$ws.Range("N13").Value = 1.421042208612402
$ws.Range("O13").Value = 3.031919051117285
$ws.Range("B14").Value = 0.5502152036388281
$ws.Range("C14").Value = 0.08783171822301483
$ws.Range("D14").Value = 0.3619087869176667
$ws.Range("F14").Value = 1.398423334518966
$ws.Range("G14").Value = 0.002442184741533182
$ws.Range("I14").Value = 0.6083420371247037
$ws.Range("J14").Value = 0.4087158797287032
$ws.Range("K14").Value = 0.6002648725875304
$ws.Range("N14").Value = 1.423739948357541
$ws.Range("O14").Value = 3.032807793595936
$ws.Range("B15").Value = 0.5451950602286502
$ws.Range("C15").Value = 0.08694305468304719
$ws.Range("D15").Value = 0.360321553761338
$ws.Range("F15").Value = 1.397743264681068
$ws.Range("G15").Value = 0.002442435678833632
$ws.Range("I15").Value = 0.6087832829720803
$ws.Range("J15").Value = 0.4069904834507554
$ws.Range("K15").Value = 0.5946320973580725
$ws.Range("N15").Value = 1.425402106383856
$ws.Range("O15").Value = 3.033384099527808
$ws.Range("B16").Value = 0.5164369894879712
$ws.Range("C16").Value = 0.08184952655875577
$ws.Range("D16").Value = 0.3512673642169659
$ws.Range("F16").Value = 1.394068471102571
$ws.Range("G16").Value = 0.002443896270069476
$ws.Range("I16").Value = 0.6114066790450998
$ws.Range("J16").Value = 0.397165175406812
$ws.Range("K16").Value = 0.5623602486084849
$ws.Range("N16").Value = 1.435076125930202
$ws.Range("O16").Value = 3.037175939053185
$ws.Range("B17").Value = 0.4988045678914261
$ws.Range("C17").Value = 0.07872401823780706
$ws.Range("D17").Value = 0.3457501506749736
$ws.Range("F17").Value = 1.392012597020781
$ws.Range("G17").Value = 0.002444812477321639
$ws.Range("I17").Value = 0.6131006104694556
$ws.Range("J17").Value = 0.3911932029319871
$ws.Range("K17").Value = 0.5425695526484446
$ws.Range("N17").Value = 1.441143621812333
$ws.Range("O17").Value = 3.039938614356373
$ws.Range("B18").Value = 0.4886666349968891
$ws.Range("C18").Value = 0.07692604052074614
$ws.Range("D18").Value = 0.3425905263914615
$ws.Range("F18").Value = 1.390903070115598
$ws.Range("G18").Value = 0.002445346889044945
$ws.Range("I18").Value = 0.6141059820991366
$ws.Range("J18").Value = 0.3877786903817224
$ws.Range("K18").Value = 0.5311892587875775
$ws.Range("N18").Value = 1.444682274010091
$ws.Range("O18").Value = 3.041688016666825
$ws.Range("B19").Value = 0.4852347711774314
$ws.Range("C19").Value = 0.07631723237940946
$ws.Range("D19").Value = 0.3415230945697658
$ws.Range("F19").Value = 1.390539935030759
$ws.Range("G19").Value = 0.002445529110037445
$ws.Range("I19").Value = 0.6144517185593656
$ws.Range("J19").Value = 0.3866261000830775
$ws.Range("K19").Value = 0.5273365855885572
$ws.Range("N19").Value = 1.445888779126705
$ws.Range("O19").Value = 3.042307877645982
$ws.Range("B20").Value = 0.5006811835540361
$ws.Range("C20").Value = 0.07905676220522651
$ws.Range("D20").Value = 0.3463360474029571
$ws.Range("F20").Value = 1.392223897769469
$ws.Range("G20").Value = 0.002444714176479448
$ws.Range("I20").Value = 0.6129170730366376
$ws.Range("J20").Value = 0.3918268174123938
$ws.Range("K20").Value = 0.5446760217722897
$ws.Range("N20").Value = 1.440492676725695
$ws.Range("O20").Value = 3.039627922987137
$ws.Range("B21").Value = 0.5526226747899159
$ws.Range("C21").Value = 0.08825783877432514
$ws.Range("D21").Value = 0.3626706426188377
$ws.Range("F21").Value = 1.398753380348353
$ws.Range("G21").Value = 0.002442064808886568
$ws.Range("I21").Value = 0.6081321249310463
$ws.Range("J21").Value = 0.4095443551494924
$ws.Range("K21").Value = 0.6029660634040397
$ws.Range("N21").Value = 1.422945532217787
$ws.Range("O21").Value = 3.032540094977236
$ws.Range("B22").Value = 0.5865988386129857
$ws.Range("C22").Value = 0.09426838377558511
$ws.Range("D22").Value = 0.3734675391591225
$ws.Range("F22").Value = 1.40367019515206
$ws.Range("G22").Value = 0.002440399134305918
$ws.Range("I22").Value = 0.6052817357110989
$ws.Range("J22").Value = 0.4213055436175495
$ws.Range("K22").Value = 0.6410824821697076
$ws.Range("N22").Value = 1.411912220761179
$ws.Range("O22").Value = 3.02933662040121
$ws.Range("B23").Value = 0.5684626327048363
$ws.Range("C23").Value = 0.09106073166469741
$ws.Range("D23").Value = 0.3676939813478555
$ws.Range("F23").Value = 1.400986589078528
$ws.Range("G23").Value = 0.002441282128883523
$ws.Range("I23").Value = 0.606777700378391
$ws.Range("J23").Value = 0.4150117446793899
$ws.Range("K23").Value = 0.6207373661922873
$ws.Range("N23").Value = 1.417761110867501
$ws.Range("O23").Value = 3.030915606674824
$ws.Range("B24").Value = 0.4998327681174999
$ws.Range("C24").Value = 0.07890633202120512
$ws.Range("D24").Value = 0.3460711251838404
$ws.Range("F24").Value = 1.39212814310558
$ws.Range("G24").Value = 0.002444758594476241
$ws.Range("I24").Value = 0.6129999521663798
$ws.Range("J24").Value = 0.3915403015808891
$ws.Range("K24").Value = 0.5437236943818107
$ws.Range("N24").Value = 1.440786812106552
$ws.Range("O24").Value = 3.039767884730168
$ws.Range("B25").Value = 0.4259506660144154
$ws.Range("C25").Value = 0.06578584816097077
$ws.Range("D25").Value = 0.323272815931773
$ws.Range("F25").Value = 1.38536224567622
$ws.Range("G25").Value = 0.002448792947562146
$ws.Range("I25").Value = 0.62089929841542
$ws.Range("J25").Value = 0.3670026579614643
$ws.Range("K25").Value = 0.4607609531287267
$ws.Range("N25").Value = 1.467488967886565
$ws.Range("O25").Value = 3.055447870799782

Write-Host "Applied 240 cell updates"
